$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the extra trailing spaces from the translated text in column B
# (rows 2 through 54 hold the language name in column A and the
# corresponding translated UI strings in column B).
for ($r = 2; $r -le 54; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val.ToString().Trim()
    }
}

# Match the author's final selection/view state
$ws.Range("B54").Select() | Out-Null

